$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 55; this shifts existing rows 55-125 down to 56-126,
# carrying formatting down from the inserted point (copies format from row above).
$ws.Rows("55:55").Insert()

# Populate the newly inserted row 55 with the new record's data
$ws.Range("A55").Value = 9
$ws.Range("B55").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C55").Value = "Metropolitana"
$ws.Range("D55").Value = (Get-Date -Year 2022 -Month 8 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E55").Value = 13
$ws.Range("F55").Value = 100112022
$ws.Range("G55").Value = "Arveja Verde"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 43
$ws.Range("K55").Value = 42000
$ws.Range("L55").Value = 42000
$ws.Range("M55").Value = 42000
$ws.Range("N55").Value = "`$/malla 25 kilos"
$ws.Range("O55").Value = "Provincia de Huasco"
$ws.Range("P55").Value = 1680
$ws.Range("Q55").Value = 25
$ws.Range("R55").Value = "Hortaliza"
